$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 10; $row -le 17; $row++) {
    $ws.Cells.Item($row, 1).Value = "Dataset Nutrition Label"
}
